$wb = $excel.ActiveWorkbook

# --- Sheet "Restricciones_del_follower" : rows 2-5, columns A:F -------------
$wsFollower = $wb.Worksheets.Item(3)
$followerRange = $wsFollower.Range("A2:F5")

# Force text storage so numeric-looking values keep their original shared
# string ("text") representation instead of being re-interpreted as numbers.
$followerRange.NumberFormat = "@"

$wsFollower.Range("A2").Value = "1.0622714025740718y"
$wsFollower.Range("B2").Value = "4.647036992919102"
$wsFollower.Range("C2").Value = "J_0_L0_v"
$wsFollower.Range("D2").Value = "0.33468162538227564"
$wsFollower.Range("E2").Value = "0"
$wsFollower.Range("F2").Value = "0.9798094181653134"

$wsFollower.Range("A3").Value = "-3 - x"
$wsFollower.Range("B3").Value = "-7.184892416399492"
$wsFollower.Range("C3").Value = "J_0_L0_v"
$wsFollower.Range("D3").Value = "0.9092567913461869"
$wsFollower.Range("E3").Value = "0.6093440204568766"
$wsFollower.Range("F3").Value = "0"

$wsFollower.Range("A4").Value = "-12 + x + 0.1552559622878642y"
$wsFollower.Range("B4").Value = "-7.135921267961507"
$wsFollower.Range("C4").Value = "J_0_LP_v"
$wsFollower.Range("D4").Value = "0.7906785535517057"
$wsFollower.Range("E4").Value = "0.8364628831988038"
$wsFollower.Range("F4").Value = "0"

$wsFollower.Range("A5").Value = "-12 + 4x + 0.34214137273234835y"
$wsFollower.Range("B5").Value = "6.236309210729873"
$wsFollower.Range("C5").Value = "J_Ne_L0_v"
$wsFollower.Range("D5").Value = "0.5618257705012442"
$wsFollower.Range("E5").Value = "0.1552176932798911"
$wsFollower.Range("F5").Value = "0.31558162870132267"

# Restore the default (General) style now that the text has been stored.
$followerRange.Style = "Normal"

# --- Sheet "Punto_modificado" : x / y sample point --------------------------
$wsPoint = $wb.Worksheets.Item(4)
$pointRange = $wsPoint.Range("A2:B2")
$pointRange.NumberFormat = "@"
$wsPoint.Range("A2").Value = "4.184892416399492"
$wsPoint.Range("B2").Value = "4.374623078112156"
$pointRange.Style = "Normal"

# --- Sheet "Vector_bf" -------------------------------------------------------
# NOTE: worksheet name lookup is case-insensitive ("Vector_bf" and
# "Vector_BF" would resolve to the same tab), so these two sheets must be
# addressed by their (1-based) tab position instead of by name.
$wsBf = $wb.Worksheets.Item(5)
$bfRange = $wsBf.Range("A2")
$bfRange.NumberFormat = "@"
$wsBf.Range("A2").Value = "-1.6705041196583514"
$bfRange.Style = "Normal"

# --- Sheet "Vector_BF" -------------------------------------------------------
$wsBF = $wb.Worksheets.Item(6)
$bFRange = $wsBF.Range("A2:A3")
$bFRange.NumberFormat = "@"
$wsBF.Range("A2").Value = "0.15201036413850844"
$wsBF.Range("A3").Value = "2.817027755499758"
$bFRange.Style = "Normal"

# --- Sheet "Vector_Alpha" ----------------------------------------------------
# This cell is a genuine numeric cell (no "t=s" in the original), so it is
# kept as a plain number.
$wsAlpha = $wb.Worksheets.Item(7)
$wsAlpha.Range("A2").Value = 0.9223720188560679
